$d = $word.ActiveDocument

# --- Helper: run a Find/Replace scoped to a given Range (in-place) ---
function Replace-InRange($range, $findText, $replaceText) {
    $range.Find.Execute($findText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replaceText, 2) | Out-Null
}

# 1. Title paragraph: "AI Tools to Generate Quizzes" -> "4. AI Tools to Generate Quizzes"
$p1 = $d.Paragraphs(1).Range
Replace-InRange $p1 "AI Tools to Generate Quizzes" "4. AI Tools to Generate Quizzes"

# 2. "Introduction" -> "4.1 Introduction"
$p2 = $d.Paragraphs(2).Range
Replace-InRange $p2 "Introduction" "4.1 Introduction"

# 3. "How AI Quiz Tools Work" -> "4.2 How AI Quiz Tools Work"
$p3 = $d.Paragraphs(3).Range
Replace-InRange $p3 "How AI Quiz Tools Work" "4.2 How AI Quiz Tools Work"

# 4. "Types of Questions Generated" -> "4.3 Types of Questions Generated"
$p4 = $d.Paragraphs(4).Range
Replace-InRange $p4 "Types of Questions Generated" "4.3 Types of Questions Generated"

# 5. "Personalization and Adaptivity" -> "4.4 Personalization and Adaptivity"
$p5 = $d.Paragraphs(5).Range
Replace-InRange $p5 "Personalization and Adaptivity" "4.4 Personalization and Adaptivity"

# 6. "Challenges and Limitations" -> "4.5 Challenges and Limitations"
$p6 = $d.Paragraphs(6).Range
Replace-InRange $p6 "Challenges and Limitations" "4.5 Challenges and Limitations"

# 7. "Future Trends" -> "4.6 Future Trends"
$p7 = $d.Paragraphs(7).Range
Replace-InRange $p7 "Future Trends" "4.6 Future Trends"
